$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.288.58"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "1.733.11"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.267"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0638"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.980.40"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "1.735.17"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "28.295.03"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "0.0₃0757"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0514"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "1.500.41"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.606"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "1.884.02"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +8.95%  "
$ws.Range("E48").Value = "  +6.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.86%  "
